# Auto-generated Excel COM-interop script to apply crypto price/volume updates
# matching the commit "Updated cryptos list on Sat Jan 20 08:47:59 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "41.524.74"
$ws.Range("E2").Value = "  +0.72%  "

# Row 3
$ws.Range("D3").Value = "2.472.60"
$ws.Range("E3").Value = "  +0.41%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.21%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.73%  "

# Row 7
$ws.Range("E7").Value = "  +0.03%  "

# Row 8
$ws.Range("E8").Value = "  -0.12%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.510"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.40%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.42"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.95%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0786"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.67%  "

# Row 12
$ws.Range("E12").Value = "  +1.12%  "

# Row 13
$ws.Range("D13").Value = "2.856.75"
$ws.Range("E13").Value = "  +0.52%  "

# Row 14
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "16.27"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.56%  "

# Row 15
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.86"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.69%  "

# Row 16
$ws.Range("D16").Value = "2.458.69"
$ws.Range("E16").Value = "  -1.00%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.769"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.15%  "

# Row 18
$ws.Range("D18").Value = "41.508.86"
$ws.Range("E18").Value = "  +0.71%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.78%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0943"
$ws.Range("E20").Value = "  +2.35%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.94%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.02"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.70%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.87"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.63%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.07%  "

# Row 25
$ws.Range("E25").Value = "  -0.04%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.05%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.60%  "

# Row 28
$ws.Range("E28").Value = "  -0.84%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.04%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.55"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.43%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "156.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.02%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.43"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.67%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.57"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.66%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0755"
$ws.Range("D34").Style = "Normal"

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.29"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.17%  "

# Row 36
$ws.Range("E36").Value = "  -8.52%  "

# Row 37
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.104"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.24%  "

# Row 38
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.88"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.72%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.81"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.89%  "

# Row 40
$ws.Range("E40").Value = "  -0.20%  "

# Row 41
$ws.Range("E41").Value = "  -4.58%  "

# Row 42
$ws.Range("E42").Value = "  -0.17%  "

# Row 43
$ws.Range("D43").Value = "1.956.93"
$ws.Range("E43").Value = "  -1.30%  "

# Row 44
$ws.Range("E44").Value = "  -0.47%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.77"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.07%  "

# Row 46
$ws.Range("E46").Value = "  -2.61%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.54%  "

# Row 48
$ws.Range("D48").Value = "2.715.42"
$ws.Range("E48").Value = "  +0.38%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.86"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.38%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.23"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.19%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.22%  "
